$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 updates
$ws.Range("I2").Value = 5.75
$ws.Range("J2").Value = 2.6
$ws.Range("M2").Value = 1.17
$ws.Range("N2").Value = 4.75
$ws.Range("X2").Value = 6.5
$ws.Range("AN2").Value = 3.4
$ws.Range("AV2").Value = 7

# Row 3 updates
$ws.Range("N3").Value = 8.1
